# Update cryptocurrency price/volume data per upstream scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.060.82'
$ws.Range('E2').Value = '  +1.00%  '

$ws.Range('D3').Value = '2.541.52'
$ws.Range('E3').Value = '  -0.22%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.21'
$ws.Range('E5').Value = '  +0.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.89'
$ws.Range('E6').Value = '  +1.71%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  -0.28%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.539'
$ws.Range('E9').Value = '  +0.67%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.06'
$ws.Range('E10').Value = '  -1.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0820'
$ws.Range('E11').Value = '  +0.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.56'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.109'
$ws.Range('E13').Value = '  -4.52%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.932.10'
$ws.Range('E14').Value = '  -0.24%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.578.01'
$ws.Range('E15').Value = '  +0.84%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.14'
$ws.Range('E16').Value = '  -3.41%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.852'
$ws.Range('E17').Value = '  -1.47%  '

$ws.Range('D18').Value = '43.085.09'
$ws.Range('E18').Value = '  +0.89%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.84'
$ws.Range('E19').Value = '  +3.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.65'
$ws.Range('E20').Value = '  -3.56%  '

$ws.Range('E21').Value = '  -0.28%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.80'
$ws.Range('E22').Value = '  -1.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.91'
$ws.Range('E23').Value = '  -0.85%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +0.18%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +1.40%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.53'
$ws.Range('E26').Value = '  -3.38%  '

$ws.Range('E27').Value = '  +0.14%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.59'
$ws.Range('E28').Value = '  +5.96%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.40'
$ws.Range('E29').Value = '  +2.10%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.40'
$ws.Range('E30').Value = '  +2.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.91'
$ws.Range('E31').Value = '  -1.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.89'
$ws.Range('E32').Value = '  +0.80%  '

$ws.Range('E33').Value = '  -0.87%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.32'
$ws.Range('E34').Value = '  -0.11%  '

$ws.Range('E35').Value = '  -1.42%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.71'
$ws.Range('E36').Value = '  +3.60%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0802'
$ws.Range('E37').Value = '  +1.38%  '

$ws.Range('E38').Value = '  +1.71%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.47'
$ws.Range('E39').Value = '  +7.31%  '

$ws.Range('E40').Value = '  -0.55%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '21.82'
$ws.Range('E41').Value = '  -9.29%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.82'
$ws.Range('E42').Value = '  -0.64%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0305'
$ws.Range('E43').Value = '  +0.67%  '

$ws.Range('E44').Value = '  +0.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  -2.13%  '

$ws.Range('D46').Value = '2.005.05'
$ws.Range('E46').Value = '  -1.79%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.12'
$ws.Range('E47').Value = '  +1.98%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '84.67'
$ws.Range('E48').Value = '  -0.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.90'
$ws.Range('E49').Value = '  +3.65%  '

$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.785.59'
$ws.Range('E50').Value = '  -0.23%  '

$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.88'
$ws.Range('E51').Value = '  +0.78%  '
